$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 21: new values in S21:V21 ---
$ws.Range("S21").Value = 6.7762000000000002
$ws.Range("T21").Value = 6.7695999999999996
$ws.Range("U21").Value = 6.7950999999999997
$ws.Range("V21").Value = 6.7763

# --- Row 22: new values in S22:V22 (T22 carries a distinct font-applied style) ---
$ws.Range("S22").Value = 6.7763
$ws.Range("T22").Value = 6.7763
$ws.Range("T22").Font.ThemeColor = 1
$ws.Range("U22").Value = 6.7760999999999996
$ws.Range("V22").Value = 6.7695999999999996

# --- Row 28: "Modified Rothermel Matlab Code" data ---
$ws.Range("B28").Value = 6
$ws.Range("C28").Value = 6
$ws.Range("D28").Value = 6
$ws.Range("E28").Value = 6
$ws.Range("F28").Value = 6
$ws.Range("G28").Value = 6
$ws.Range("H28").Value = 6
$ws.Range("I28").Value = 6
$ws.Range("J28").Value = 6
$ws.Range("K28").Value = 6
$ws.Range("L28").Value = 4.8231000000000002
$ws.Range("M28").Value = 6
$ws.Range("N28").Value = 4.8121
$ws.Range("O28").Value = 6
$ws.Range("P28").Value = 6
$ws.Range("Q28").Value = 6
$ws.Range("R28").Value = 6
$ws.Range("S28").Value = 6
$ws.Range("T28").Value = 6
$ws.Range("U28").Value = 6
$ws.Range("V28").Value = 6

# --- Row 29: "Original Rothermel Matlab Code" data ---
$ws.Range("B29").Value = 6
$ws.Range("C29").Value = 6
$ws.Range("D29").Value = 6
$ws.Range("E29").Value = 6
$ws.Range("F29").Value = 6
$ws.Range("G29").Value = 6
$ws.Range("H29").Value = 6
$ws.Range("I29").Value = 6
$ws.Range("J29").Value = 6
$ws.Range("K29").Value = 6
$ws.Range("L29").Value = 6
$ws.Range("M29").Value = 6
$ws.Range("N29").Value = 6
$ws.Range("O29").Value = 6
$ws.Range("P29").Value = 6
$ws.Range("Q29").Value = 6
$ws.Range("R29").Value = 6
$ws.Range("S29").Value = 6
$ws.Range("T29").Value = 6
$ws.Range("U29").Value = 6
$ws.Range("V29").Value = 6

# --- Row 30: "ROS var modified" data ---
$ws.Range("B30").Value = 6.4665999999999997
$ws.Range("C30").Value = 6.5193000000000003
$ws.Range("E30").Value = 6.3959999999999999
$ws.Range("G30").Value = 6.5202999999999998
$ws.Range("I30").Value = 6.3602999999999996
$ws.Range("K30").Value = 6.4885000000000002
$ws.Range("L30").Value = 4.8231000000000002
$ws.Range("M30").Value = 6.4313000000000002
$ws.Range("N30").Value = 4.8121
$ws.Range("O30").Value = 6.4466999999999999
$ws.Range("P30").Value = 8.0684000000000005
$ws.Range("Q30").Value = 6.5022000000000002
$ws.Range("R30").Value = 8.0904000000000007
$ws.Range("S30").Value = 6.4667000000000003
$ws.Range("T30").Value = 6.4608999999999996
$ws.Range("U30").Value = 6.4664999999999999
$ws.Range("V30").Value = 6.4607999999999999

# --- Row 31: "ROS var original" data (B31 already existed w/ Menlo style s="1") ---
$ws.Range("B31").Value = 6.4665999999999997
$ws.Range("C31").Value = 6.5358999999999998
$ws.Range("E31").Value = 6.4668999999999999
$ws.Range("G31").Value = 6.5453000000000001
$ws.Range("I31").Value = 6.4664999999999999
$ws.Range("K31").Value = 6.4968000000000004
$ws.Range("L31").Value = 6.3724999999999996
$ws.Range("M31").Value = 6.4667000000000003
$ws.Range("N31").Value = 6.4161000000000001
$ws.Range("O31").Value = 6.4383999999999997
$ws.Range("P31").Value = 6.4897
$ws.Range("Q31").Value = 6.4668000000000001
$ws.Range("R31").Value = 6.4722999999999997
$ws.Range("S31").Value = 6.4668000000000001
$ws.Range("T31").Value = 6.4672000000000001
$ws.Range("U31").Value = 6.4665999999999997
$ws.Range("V31").Value = 6.4672999999999998

# --- View state: zoom to 100%, scroll back to top-left, select V30 ---
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$win.Zoom = 100
$ws.Range("V30").Select()
